# Automatic BRVM update (via GitHub Actions) applied through Excel COM interop.
# Updates both worksheets: "Recommandations" and "Top_YTD".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update index rows 2-23 (D/E values only) ---
$ws1.Cells.Item(2,4).Value = 3925
$ws1.Cells.Item(2,5).Value = 985

$ws1.Cells.Item(3,4).Value = 3338.78
$ws1.Cells.Item(3,5).Value = 111.44

$ws1.Cells.Item(4,4).Value = 2775

$ws1.Cells.Item(5,4).Value = 2690

$ws1.Cells.Item(6,4).Value = 2643.7
$ws1.Cells.Item(6,5).Value = 653.39

$ws1.Cells.Item(7,4).Value = 2370

$ws1.Cells.Item(8,4).Value = 2330
$ws1.Cells.Item(8,5).Value = 575

$ws1.Cells.Item(9,4).Value = 2225

$ws1.Cells.Item(10,4).Value = 2125

$ws1.Cells.Item(11,4).Value = 1476.01
$ws1.Cells.Item(11,5).Value = 363.76

$ws1.Cells.Item(12,4).Value = 1392.76
$ws1.Cells.Item(12,5).Value = 341.46

$ws1.Cells.Item(13,4).Value = 1287.41
$ws1.Cells.Item(13,5).Value = 315.73

$ws1.Cells.Item(14,4).Value = 801.62
$ws1.Cells.Item(14,5).Value = 199.06

$ws1.Cells.Item(15,4).Value = 701.76
$ws1.Cells.Item(15,5).Value = 174.55

$ws1.Cells.Item(16,4).Value = 698.17
$ws1.Cells.Item(16,5).Value = 172.59

$ws1.Cells.Item(17,4).Value = 531.7
$ws1.Cells.Item(17,5).Value = 131.27

$ws1.Cells.Item(18,4).Value = 520.26
$ws1.Cells.Item(18,5).Value = 130.63

$ws1.Cells.Item(19,4).Value = 485
$ws1.Cells.Item(19,5).Value = 121.75

$ws1.Cells.Item(20,4).Value = 476.66
$ws1.Cells.Item(20,5).Value = 119.66

$ws1.Cells.Item(21,4).Value = 438.15
$ws1.Cells.Item(21,5).Value = 107.98

$ws1.Cells.Item(22,4).Value = 423.29
$ws1.Cells.Item(22,5).Value = 104.97

$ws1.Cells.Item(23,4).Value = 376.33
$ws1.Cells.Item(23,5).Value = 94.16

# --- Row 25: TRACTAFRIC MOTORS CI (PRSC) - only C, D, G change ---
$ws1.Cells.Item(25,3).Value = 0
$ws1.Cells.Item(25,4).Value = 10
$ws1.Cells.Item(25,7).Value = "➖ Neutre"

# --- Row 26: BANK OF AFRICA NG (BOAN) - only B, D, E change ---
$ws1.Cells.Item(26,2).Value = 2
$ws1.Cells.Item(26,4).Value = 9.23
$ws1.Cells.Item(26,5).Value = 3.4

# --- Rows 27-48: full reshuffle (sorted by Variation Totale desc), row 49 removed ---
$ws1.Cells.Item(27,1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(27,2).Value = 1
$ws1.Cells.Item(27,3).Value = 0
$ws1.Cells.Item(27,4).Value = 3.73
$ws1.Cells.Item(27,5).Value = 3.73
$ws1.Cells.Item(27,6).Value = "🟡 Observer"
$ws1.Cells.Item(27,7).Value = "➖ Neutre"

$ws1.Cells.Item(28,1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Cells.Item(28,2).Value = 1
$ws1.Cells.Item(28,3).Value = 0
$ws1.Cells.Item(28,4).Value = 3.3
$ws1.Cells.Item(28,5).Value = 3.3
$ws1.Cells.Item(28,6).Value = "🟡 Observer"
$ws1.Cells.Item(28,7).Value = "➖ Neutre"

$ws1.Cells.Item(29,1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(29,2).Value = 1
$ws1.Cells.Item(29,3).Value = 0
$ws1.Cells.Item(29,4).Value = 2.59
$ws1.Cells.Item(29,5).Value = 2.59
$ws1.Cells.Item(29,6).Value = "🟡 Observer"
$ws1.Cells.Item(29,7).Value = "➖ Neutre"

$ws1.Cells.Item(30,1).Value = "SONATEL SN (SNTS)"
$ws1.Cells.Item(30,2).Value = 1
$ws1.Cells.Item(30,3).Value = 1
$ws1.Cells.Item(30,4).Value = 2.17
$ws1.Cells.Item(30,5).Value = 4
$ws1.Cells.Item(30,6).Value = "🟡 Observer"
$ws1.Cells.Item(30,7).Value = "👀 À surveiller"

$ws1.Cells.Item(31,1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Cells.Item(31,2).Value = 2
$ws1.Cells.Item(31,3).Value = 1
$ws1.Cells.Item(31,4).Value = 1.41
$ws1.Cells.Item(31,5).Value = -2.5
$ws1.Cells.Item(31,6).Value = "🟡 Observer"
$ws1.Cells.Item(31,7).Value = "👀 À surveiller"

$ws1.Cells.Item(32,1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(32,2).Value = 1
$ws1.Cells.Item(32,3).Value = 0
$ws1.Cells.Item(32,4).Value = 0.88
$ws1.Cells.Item(32,5).Value = 0.88
$ws1.Cells.Item(32,6).Value = "🟡 Observer"
$ws1.Cells.Item(32,7).Value = "➖ Neutre"

$ws1.Cells.Item(33,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(33,2).Value = 1
$ws1.Cells.Item(33,3).Value = 1
$ws1.Cells.Item(33,4).Value = 0.51
$ws1.Cells.Item(33,5).Value = 7.41
$ws1.Cells.Item(33,6).Value = "🟡 Observer"
$ws1.Cells.Item(33,7).Value = "👀 À surveiller"

$ws1.Cells.Item(34,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(34,2).Value = 1
$ws1.Cells.Item(34,3).Value = 1
$ws1.Cells.Item(34,4).Value = 0.28
$ws1.Cells.Item(34,5).Value = 5.36
$ws1.Cells.Item(34,6).Value = "🟡 Observer"
$ws1.Cells.Item(34,7).Value = "👀 À surveiller"

$ws1.Cells.Item(35,1).Value = "SICABLE CI (CABC)"
$ws1.Cells.Item(35,2).Value = 1
$ws1.Cells.Item(35,3).Value = 1
$ws1.Cells.Item(35,4).Value = 0.27
$ws1.Cells.Item(35,5).Value = 5.31
$ws1.Cells.Item(35,6).Value = "🟡 Observer"
$ws1.Cells.Item(35,7).Value = "👀 À surveiller"

$ws1.Cells.Item(36,1).Value = "TOTAL"
$ws1.Cells.Item(36,2).Value = 0
$ws1.Cells.Item(36,3).Value = 4
$ws1.Cells.Item(36,4).Value = 0
$ws1.Cells.Item(36,5).Value = 0
$ws1.Cells.Item(36,6).Value = "🟡 Observer"
$ws1.Cells.Item(36,7).Value = "➖ Neutre"

$ws1.Cells.Item(37,1).Value = "SODE CI (SDCC)"
$ws1.Cells.Item(37,2).Value = 1
$ws1.Cells.Item(37,3).Value = 1
$ws1.Cells.Item(37,4).Value = -0.6899999999999999
$ws1.Cells.Item(37,5).Value = 3.36
$ws1.Cells.Item(37,6).Value = "🟡 Observer"
$ws1.Cells.Item(37,7).Value = "👀 À surveiller"

$ws1.Cells.Item(38,1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(38,2).Value = 1
$ws1.Cells.Item(38,3).Value = 1
$ws1.Cells.Item(38,4).Value = -0.71
$ws1.Cells.Item(38,5).Value = -2.11
$ws1.Cells.Item(38,6).Value = "🟡 Observer"
$ws1.Cells.Item(38,7).Value = "👀 À surveiller"

$ws1.Cells.Item(39,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(39,2).Value = 1
$ws1.Cells.Item(39,3).Value = 2
$ws1.Cells.Item(39,4).Value = -0.83
$ws1.Cells.Item(39,5).Value = -2.35
$ws1.Cells.Item(39,6).Value = "🟡 Observer"
$ws1.Cells.Item(39,7).Value = "👀 À surveiller"

$ws1.Cells.Item(40,1).Value = "BICI CI (BICC)"
$ws1.Cells.Item(40,2).Value = 0
$ws1.Cells.Item(40,3).Value = 1
$ws1.Cells.Item(40,4).Value = -2.54
$ws1.Cells.Item(40,5).Value = -2.54
$ws1.Cells.Item(40,6).Value = "🟡 Observer"
$ws1.Cells.Item(40,7).Value = "➖ Neutre"

$ws1.Cells.Item(41,1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(41,2).Value = 0
$ws1.Cells.Item(41,3).Value = 1
$ws1.Cells.Item(41,4).Value = -2.86
$ws1.Cells.Item(41,5).Value = -2.86
$ws1.Cells.Item(41,6).Value = "🟡 Observer"
$ws1.Cells.Item(41,7).Value = "➖ Neutre"

$ws1.Cells.Item(42,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(42,2).Value = 0
$ws1.Cells.Item(42,3).Value = 1
$ws1.Cells.Item(42,4).Value = -3.38
$ws1.Cells.Item(42,5).Value = -3.38
$ws1.Cells.Item(42,6).Value = "🟡 Observer"
$ws1.Cells.Item(42,7).Value = "➖ Neutre"

$ws1.Cells.Item(43,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(43,2).Value = 1
$ws1.Cells.Item(43,3).Value = 2
$ws1.Cells.Item(43,4).Value = -3.77
$ws1.Cells.Item(43,5).Value = -2.69
$ws1.Cells.Item(43,6).Value = "🟡 Observer"
$ws1.Cells.Item(43,7).Value = "👀 À surveiller"

$ws1.Cells.Item(44,1).Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Cells.Item(44,2).Value = 0
$ws1.Cells.Item(44,3).Value = 1
$ws1.Cells.Item(44,4).Value = -5.69
$ws1.Cells.Item(44,5).Value = -5.69
$ws1.Cells.Item(44,6).Value = "🟡 Observer"
$ws1.Cells.Item(44,7).Value = "➖ Neutre"

$ws1.Cells.Item(45,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(45,2).Value = 0
$ws1.Cells.Item(45,3).Value = 1
$ws1.Cells.Item(45,4).Value = -5.88
$ws1.Cells.Item(45,5).Value = -5.88
$ws1.Cells.Item(45,6).Value = "🟡 Observer"
$ws1.Cells.Item(45,7).Value = "➖ Neutre"

$ws1.Cells.Item(46,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(46,2).Value = 0
$ws1.Cells.Item(46,3).Value = 1
$ws1.Cells.Item(46,4).Value = -6.09
$ws1.Cells.Item(46,5).Value = -6.09
$ws1.Cells.Item(46,6).Value = "🟡 Observer"
$ws1.Cells.Item(46,7).Value = "➖ Neutre"

$ws1.Cells.Item(47,1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(47,2).Value = 0
$ws1.Cells.Item(47,3).Value = 1
$ws1.Cells.Item(47,4).Value = -7.46
$ws1.Cells.Item(47,5).Value = -7.46
$ws1.Cells.Item(47,6).Value = "🟡 Observer"
$ws1.Cells.Item(47,7).Value = "➖ Neutre"

$ws1.Cells.Item(48,1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(48,2).Value = 0
$ws1.Cells.Item(48,3).Value = 2
$ws1.Cells.Item(48,4).Value = -12.35
$ws1.Cells.Item(48,5).Value = -4.87
$ws1.Cells.Item(48,6).Value = "🟡 Observer"
$ws1.Cells.Item(48,7).Value = "➖ Neutre"

# --- Sheet "Top_YTD": update B column values ---
$ws2.Cells.Item(2,2).Value = 9198694.99
$ws2.Cells.Item(3,2).Value = 1366599.36
$ws2.Cells.Item(4,2).Value = 396700
$ws2.Cells.Item(5,2).Value = 355995.55
$ws2.Cells.Item(6,2).Value = 335120.02
$ws2.Cells.Item(7,2).Value = 229819.9
$ws2.Cells.Item(8,2).Value = 216823.4
$ws2.Cells.Item(9,2).Value = 185141.6
$ws2.Cells.Item(10,2).Value = 158679.84
$ws2.Cells.Item(11,2).Value = 48277.32

# Row 49 (SAPH CI (SPHC)) no longer present in the refreshed list:
# it has been folded into the resorted rows 27-48 above, so the
# now-superfluous trailing row is removed and the sheet dimension
# shrinks from A1:G49 to A1:G48 automatically.
$ws1.Rows.Item(49).Delete() | Out-Null
